# Auto-generated edit script applying the Gilgamesh_Profits.xlsx market-data refresh
# (scheduled runner update) described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Making Ends Meet | Superior Spiritbond Potion (row 112)
$ws.Cells.Item(112, 8).Value = 2402.95   # H112 was 2435.7896
$ws.Cells.Item(112, 10).Value = 2402.95   # J112 was 2435.7896
$ws.Cells.Item(112, 12).Value = 7208.849999999999   # L112 was 7307.3688
$ws.Cells.Item(112, 14).Value = -9424.849999999999   # N112 was -9523.3688

# Cutting Edge of Culinary Quality | Magnesia Whetstone (row 137)
$ws.Cells.Item(137, 8).Value = 6264060.5   # H137 was 5568165
$ws.Cells.Item(137, 9).Value = 10001096   # I137 was 8334413.5
$ws.Cells.Item(137, 11).Value = 30003288   # K137 was 25003240.5
$ws.Cells.Item(137, 13).Value = -30000738   # M137 was -25000690.5

# All-night Crafting | Cunning Craftsman's Tisane (row 138)
$ws.Cells.Item(138, 8).Value = 354912.84   # H138 was 345289.7
$ws.Cells.Item(138, 9).Value = 2639.394   # I138 was 2783.6453
$ws.Cells.Item(138, 10).Value = 660834.5   # J138 was 598091.75
$ws.Cells.Item(138, 11).Value = 7918.181999999999   # K138 was 8350.9359
$ws.Cells.Item(138, 12).Value = 1982503.5   # L138 was 1794275.25
$ws.Cells.Item(138, 13).Value = -2778.181999999999   # M138 was -3210.9359
$ws.Cells.Item(138, 14).Value = -1992783.5   # N138 was -1804555.25

# Something Salty and Ceremonial | Gomphotherium Codex (row 139)
$ws.Cells.Item(139, 8).Value = 69953.63   # H139 was 69998.17999999999
$ws.Cells.Item(139, 10).Value = 69953.63   # J139 was 69998.17999999999
$ws.Cells.Item(139, 12).Value = 69953.63   # L139 was 69998.17999999999
$ws.Cells.Item(139, 14).Value = -80233.63   # N139 was -80278.17999999999

# Tome for Tradition | Book of Ra'Kaznar (row 140)
$ws.Cells.Item(140, 8).Value = 59999   # H140 was 74852.71000000001
$ws.Cells.Item(140, 10).Value = 59999   # J140 was 74852.71000000001
$ws.Cells.Item(140, 12).Value = 59999   # L140 was 74852.71000000001
$ws.Cells.Item(140, 14).Value = -70359   # N140 was -85212.71000000001

$ws = $wb.Worksheets.Item("ARM")
# Dealing with the Tough Stuff | Cobalt Ingot (row 61)
$ws.Cells.Item(61, 8).Value = 4799.6   # H61 was 4564.625
$ws.Cells.Item(61, 9).Value = 2798   # I61 was 2638.182
$ws.Cells.Item(61, 11).Value = 2798   # K61 was 2638.182
$ws.Cells.Item(61, 13).Value = -2586   # M61 was -2426.182

# As the Bolt Flies | Titanium Nugget (row 74)
$ws.Cells.Item(74, 8).Value = 224152.64   # H74 was 207734.03
$ws.Cells.Item(74, 10).Value = 3845.5   # J74 was 3397.5
$ws.Cells.Item(74, 12).Value = 3845.5   # L74 was 3397.5
$ws.Cells.Item(74, 14).Value = -5593.5   # N74 was -5145.5

# Heavy Metal Banned (L) | Titanium Nugget (row 77)
$ws.Cells.Item(77, 8).Value = 224152.64   # H77 was 207734.03
$ws.Cells.Item(77, 10).Value = 3845.5   # J77 was 3397.5
$ws.Cells.Item(77, 12).Value = 19227.5   # L77 was 16987.5
$ws.Cells.Item(77, 14).Value = -27963.5   # N77 was -25723.5

# Signed, Shield, Delivered | Titanbronze Tower Shield (row 117)
$ws.Cells.Item(117, 8).Value = 40000   # H117 was 39247.5
$ws.Cells.Item(117, 10).Value = 40000   # J117 was 39247.5
$ws.Cells.Item(117, 12).Value = 40000   # L117 was 39247.5
$ws.Cells.Item(117, 14).Value = -49178   # N117 was -48425.5

# Haste for High Durium | High Durium Nugget (row 122)
$ws.Cells.Item(122, 8).Value = 2763.743   # H122 was 2593.1282
$ws.Cells.Item(122, 9).Value = 2550.9688   # I122 was 2389.7778
$ws.Cells.Item(122, 11).Value = 7652.9064   # K122 was 7169.3334
$ws.Cells.Item(122, 13).Value = -5202.9064   # M122 was -4719.3334

# Don't Bore Me, Ore Me | Mountain Chromite Ingot (row 132)
$ws.Cells.Item(132, 8).Value = 2900.1924   # H132 was 2673.0334
$ws.Cells.Item(132, 9).Value = 1837.7333   # I132 was 1702.7368
$ws.Cells.Item(132, 11).Value = 5513.199900000001   # K132 was 5108.2104
$ws.Cells.Item(132, 13).Value = -2983.199900000001   # M132 was -2578.2104

# Shielding My Students | Mountain Chromite Tower Shield (row 133)
$ws.Cells.Item(133, 8).Value = 99487   # H133 was 99494.5
$ws.Cells.Item(133, 10).Value = 99487   # J133 was 99494.5
$ws.Cells.Item(133, 12).Value = 99487   # L133 was 99494.5
$ws.Cells.Item(133, 14).Value = -104547   # N133 was -104554.5

# Brace for More Vambraces | Ruthenium Vambraces of Maiming (row 134)
$ws.Cells.Item(134, 8).Value = 101995.25   # H134 was 103999
$ws.Cells.Item(134, 10).Value = 101995.25   # J134 was 103999
$ws.Cells.Item(134, 12).Value = 101995.25   # L134 was 103999
$ws.Cells.Item(134, 14).Value = -112135.25   # N134 was -114139

# Forgiveness for My Shins | Ruthenium Sabatons of Fending (row 135)
$ws.Cells.Item(135, 8).Value = 107498.5   # H135 was 100000
$ws.Cells.Item(135, 10).Value = 107498.5   # J135 was 100000
$ws.Cells.Item(135, 12).Value = 107498.5   # L135 was 100000
$ws.Cells.Item(135, 14).Value = -117638.5   # N135 was -110140

# Metal with Mettle | Cobalt Tungsten Ingot (row 136)
$ws.Cells.Item(136, 8).Value = 4799.6   # H136 was 4564.625
$ws.Cells.Item(136, 9).Value = 2798   # I136 was 2638.182
$ws.Cells.Item(136, 11).Value = 8394   # K136 was 7914.545999999999
$ws.Cells.Item(136, 13).Value = -5844   # M136 was -5364.545999999999

# Odd Instruments | Cobalt Tungsten Alembic (row 137)
$ws.Cells.Item(137, 8).Value = 100000   # H137 was 135000
$ws.Cells.Item(137, 10).Value = 100000   # J137 was 135000
$ws.Cells.Item(137, 12).Value = 100000   # L137 was 135000
$ws.Cells.Item(137, 14).Value = -110200   # N137 was -145200

# Don't Ask about the Rivets | Titanium Gold Helm of Casting (row 138)
$ws.Cells.Item(138, 8).Value = 90666.664   # H138 was 96661.336
$ws.Cells.Item(138, 10).Value = 90666.664   # J138 was 96661.336
$ws.Cells.Item(138, 12).Value = 90666.664   # L138 was 96661.336
$ws.Cells.Item(138, 14).Value = -100946.664   # N138 was -106941.336

# A Hand for a Deckhand | Ra'Kaznar Gloves of Scouting (row 140)
$ws.Cells.Item(140, 8).Value = 77500   # H140 was 80000
$ws.Cells.Item(140, 10).Value = 77500   # J140 was 80000
$ws.Cells.Item(140, 12).Value = 77500   # L140 was 80000
$ws.Cells.Item(140, 14).Value = -87860   # N140 was -90360

$ws = $wb.Worksheets.Item("CRP")
# Wall Not Found | Walnut Lumber (row 31)
$ws.Cells.Item(31, 8).Value = 3398.0425   # H31 was 3457.4375
$ws.Cells.Item(31, 9).Value = 2413.3784   # I31 was 2467.9167
$ws.Cells.Item(31, 10).Value = 7041.3   # J31 was 6426
$ws.Cells.Item(31, 11).Value = 2413.3784   # K31 was 2467.9167
$ws.Cells.Item(31, 12).Value = 7041.3   # L31 was 6426
$ws.Cells.Item(31, 13).Value = -2118.3784   # M31 was -2172.9167
$ws.Cells.Item(31, 14).Value = -7631.3   # N31 was -7016

# Armoires of the Rich and Famous | Walnut Lumber (row 34)
$ws.Cells.Item(34, 8).Value = 3398.0425   # H34 was 3457.4375
$ws.Cells.Item(34, 9).Value = 2413.3784   # I34 was 2467.9167
$ws.Cells.Item(34, 10).Value = 7041.3   # J34 was 6426
$ws.Cells.Item(34, 11).Value = 2413.3784   # K34 was 2467.9167
$ws.Cells.Item(34, 12).Value = 7041.3   # L34 was 6426
$ws.Cells.Item(34, 13).Value = -2211.3784   # M34 was -2265.9167
$ws.Cells.Item(34, 14).Value = -7445.3   # N34 was -6830

# You Do the Heavy Lifting | Mahogany Lumber (row 58)
$ws.Cells.Item(58, 8).Value = 2884.4827   # H58 was 2900.0688
$ws.Cells.Item(58, 9).Value = 2532.7273   # I58 was 2573.6365
$ws.Cells.Item(58, 10).Value = 3099.4443   # J58 was 3099.5557
$ws.Cells.Item(58, 11).Value = 2532.7273   # K58 was 2573.6365
$ws.Cells.Item(58, 12).Value = 3099.4443   # L58 was 3099.5557
$ws.Cells.Item(58, 13).Value = -2329.7273   # M58 was -2370.6365
$ws.Cells.Item(58, 14).Value = -3505.4443   # N58 was -3505.5557

# Wood You Be Quiet | Ceiba Lumber (row 134)
$ws.Cells.Item(134, 8).Value = 2264.25   # H134 was 2228.7222
$ws.Cells.Item(134, 9).Value = 2213.3794   # I134 was 2156.1667
$ws.Cells.Item(134, 10).Value = 2475   # J134 was 2591.5
$ws.Cells.Item(134, 11).Value = 6640.138199999999   # K134 was 6468.500100000001
$ws.Cells.Item(134, 12).Value = 7425   # L134 was 7774.5
$ws.Cells.Item(134, 13).Value = -4105.138199999999   # M134 was -3933.500100000001
$ws.Cells.Item(134, 14).Value = -12495   # N134 was -12844.5

# Turali Quality | Dark Mahogany Lumber (row 136)
$ws.Cells.Item(136, 8).Value = 2884.4827   # H136 was 2900.0688
$ws.Cells.Item(136, 9).Value = 2532.7273   # I136 was 2573.6365
$ws.Cells.Item(136, 10).Value = 3099.4443   # J136 was 3099.5557
$ws.Cells.Item(136, 11).Value = 7598.1819   # K136 was 7720.9095
$ws.Cells.Item(136, 12).Value = 9298.332900000001   # L136 was 9298.667099999999
$ws.Cells.Item(136, 13).Value = -5048.1819   # M136 was -5170.9095
$ws.Cells.Item(136, 14).Value = -14398.3329   # N136 was -14398.6671

$ws = $wb.Worksheets.Item("GSM")
# Copious Crystal Cannons | Manasilver Nugget (row 113)
$ws.Cells.Item(113, 8).Value = 5059.25   # H113 was 4243.6665
$ws.Cells.Item(113, 9).Value = 4829.6665   # I113 was 3157.3333
$ws.Cells.Item(113, 10).Value = 5197   # J113 was 5330
$ws.Cells.Item(113, 11).Value = 4829.6665   # K113 was 3157.3333
$ws.Cells.Item(113, 12).Value = 5197   # L113 was 5330
$ws.Cells.Item(113, 13).Value = -2659.6665   # M113 was -987.3332999999998
$ws.Cells.Item(113, 14).Value = -9537   # N113 was -9670

$ws = $wb.Worksheets.Item("LTW")
# Skin off Their Backs | Aldgoat Leather (row 22)
$ws.Cells.Item(22, 8).Value = 1150   # H22 was 1208.3334
$ws.Cells.Item(22, 10).Value = 1250   # J22 was 1300
$ws.Cells.Item(22, 12).Value = 1250   # L22 was 1300
$ws.Cells.Item(22, 14).Value = -1840   # N22 was -1890

# Fire and Hide | Aldgoat Leather (row 27)
$ws.Cells.Item(27, 8).Value = 1150   # H27 was 1208.3334
$ws.Cells.Item(27, 10).Value = 1250   # J27 was 1300
$ws.Cells.Item(27, 12).Value = 1250   # L27 was 1300
$ws.Cells.Item(27, 14).Value = -1464   # N27 was -1514

# Best Served Toad | Toad Leather (row 40)
$ws.Cells.Item(40, 8).Value = 61549.723   # H40 was 55996.55
$ws.Cells.Item(40, 9).Value = 68118.75   # I40 was 61218.668
$ws.Cells.Item(40, 11).Value = 68118.75   # K40 was 61218.668
$ws.Cells.Item(40, 13).Value = -67982.75   # M40 was -61082.668

# Spelling Me Softly | Raptor Leather (row 61)
$ws.Cells.Item(61, 8).Value = 4999   # H61 was 2254.55
$ws.Cells.Item(61, 9).Value = 4999   # I61 was 2127
$ws.Cells.Item(61, 10).Value = 0   # J61 was 3402.5
$ws.Cells.Item(61, 11).Value = 4999   # K61 was 2127
$ws.Cells.Item(61, 12).Value = 0   # L61 was 3402.5
$ws.Cells.Item(61, 14).Value = -4797   # N61 was -3806.5

# Peace in Rest | Atrociraptor Leather (row 113)
$ws.Cells.Item(113, 8).Value = 4999   # H113 was 2254.55
$ws.Cells.Item(113, 9).Value = 4999   # I113 was 2127
$ws.Cells.Item(113, 10).Value = 0   # J113 was 3402.5
$ws.Cells.Item(113, 11).Value = 4999   # K113 was 2127
$ws.Cells.Item(113, 12).Value = 0   # L113 was 3402.5
$ws.Cells.Item(113, 14).Value = -2829   # N113 was -7742.5

# Hell on Leather | Gaja Leather (row 122)
$ws.Cells.Item(122, 8).Value = 0   # H122 was 5999
$ws.Cells.Item(122, 9).Value = 0   # I122 was 5999
$ws.Cells.Item(122, 11).Value = 0   # K122 was 17997

# Tenets of Tanning | Silver Lobo Leather (row 132)
$ws.Cells.Item(132, 8).Value = 6246.9414   # H132 was 6258.7646
$ws.Cells.Item(132, 9).Value = 3313.8572   # I132 was 3534.5
$ws.Cells.Item(132, 10).Value = 8300.1   # J132 was 7744.727
$ws.Cells.Item(132, 11).Value = 9941.571599999999   # K132 was 10603.5
$ws.Cells.Item(132, 12).Value = 24900.3   # L132 was 23234.181
$ws.Cells.Item(132, 13).Value = -7411.571599999999   # M132 was -8073.5
$ws.Cells.Item(132, 14).Value = -29960.3   # N132 was -28294.181

# Respect for Br'aax | Br'aax Leather (row 136)
$ws.Cells.Item(136, 8).Value = 4933.5625   # H136 was 4946.125
$ws.Cells.Item(136, 9).Value = 5798.778   # I136 was 5821.1113
$ws.Cells.Item(136, 11).Value = 17396.334   # K136 was 17463.3339
$ws.Cells.Item(136, 13).Value = -14846.334   # M136 was -14913.3339

$ws = $wb.Worksheets.Item("WVR")
# Pride Up in Smoke | Rainbow Cloth (row 62)
$ws.Cells.Item(62, 8).Value = 10499.2   # H62 was 11124.25
$ws.Cells.Item(62, 9).Value = 8166.3335   # I62 was 8250
$ws.Cells.Item(62, 11).Value = 8166.3335   # K62 was 8250
$ws.Cells.Item(62, 13).Value = -7542.3335   # M62 was -7626

# Desperate for Diversionaries (L) | Rainbow Cloth (row 65)
$ws.Cells.Item(65, 8).Value = 10499.2   # H65 was 11124.25
$ws.Cells.Item(65, 9).Value = 8166.3335   # I65 was 8250
$ws.Cells.Item(65, 11).Value = 40831.6675   # K65 was 41250
$ws.Cells.Item(65, 13).Value = -37711.6675   # M65 was -38130

# A Tender Table | Pixie Floss (row 113)
$ws.Cells.Item(113, 8).Value = 939   # H113 was 860.7692
$ws.Cells.Item(113, 9).Value = 900   # I113 was 739
$ws.Cells.Item(113, 10).Value = 948.75   # J113 was 936.875
$ws.Cells.Item(113, 11).Value = 2700   # K113 was 2217
$ws.Cells.Item(113, 12).Value = 2846.25   # L113 was 2810.625
$ws.Cells.Item(113, 13).Value = -530   # M113 was -47
$ws.Cells.Item(113, 14).Value = -7186.25   # N113 was -7150.625

# Comfy Cabins | Snow Cotton Cloth (row 132)
$ws.Cells.Item(132, 8).Value = 3726.606   # H132 was 4282.0713
$ws.Cells.Item(132, 9).Value = 3569.4827   # I132 was 4184.7915
$ws.Cells.Item(132, 11).Value = 10708.4481   # K132 was 12554.3745
$ws.Cells.Item(132, 13).Value = -8178.4481   # M132 was -10024.3745

# Weaving the Envelope | Sarcenet Cloth (row 136)
$ws.Cells.Item(136, 8).Value = 43481744   # H136 was 43481750
$ws.Cells.Item(136, 9).Value = 55556670   # I136 was 55556684
$ws.Cells.Item(136, 11).Value = 166670010   # K136 was 166670052
$ws.Cells.Item(136, 13).Value = -166667460   # M136 was -166667502
